# "Fixed the goddamn bug" - add the two missing client rows to the
# Clients sheet (rows 3 and 4), matching the Id/Balance/Tariff pattern
# already present in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clients")

# Row 3: new client "STACK"
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "STACK"
$ws.Cells.Item(3, 3).Value = "STACK"
$ws.Cells.Item(3, 4).Value = "STACK"
$ws.Cells.Item(3, 5).Value = 100
$ws.Cells.Item(3, 6).Value = 101
$ws.Cells.Item(3, 7).Value = 555

# Row 4: new client "TEAN" / "TEAM"
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "TEAN"
$ws.Cells.Item(4, 3).Value = "TEAM"
$ws.Cells.Item(4, 4).Value = "TEAM"
$ws.Cells.Item(4, 5).Value = 100
$ws.Cells.Item(4, 6).Value = 101
$ws.Cells.Item(4, 7).Value = 555

# Leave the selection where the author ended up after entering the data.
$ws.Range("A8").Select()
